$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 86; this shifts existing rows 86..122 down to 87..123
# and extends the used range to A1:R123 automatically.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with a duplicate of the record that
# used to be on row 86 (now shifted to row 87), but with a new date.
$ws.Range("A86").Value = 7
$ws.Range("B86").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C86").Value = "Ñuble"
$ws.Range("D86").Value = 44460
$ws.Range("E86").Value = 16
$ws.Range("F86").Value = 100112006
$ws.Range("G86").Value = "Repollo"
$ws.Range("H86").Value = "Crespo record"
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 300
$ws.Range("K86").Value = 600
$ws.Range("L86").Value = 650
$ws.Range("M86").Value = 625
$ws.Range("N86").Value = "$/unidad"
$ws.Range("O86").Value = "Provincia de Diguillín"
$ws.Range("P86").Value = 625
$ws.Range("Q86").Value = 1
$ws.Range("R86").Value = "Hortaliza"

# Match the date-cell number format used by the rest of column D.
$ws.Range("D86").NumberFormat = $ws.Range("D87").NumberFormat
